$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-24 10:45:23"
$wsZh.Range("H2").Value = "2016-03-24 10:46:04"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-24 10:45:33"
$wsDe.Range("H2").Value = "2016-03-24 10:46:18"
